# "update lector de .txt" — the code/description list now comes from a .txt
# reader that emits zero-padded numeric IDs (00000001, 00000002) instead of
# the old part codes (MK-5N, MK-7N). Re-style those ID cells like the
# generated output (Roboto / #333333, stored as text so the leading zeros
# survive) and leave everything else (headers, quantities, descriptions)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new "text id" format once on a scratch cell, then copy that
# format onto A2/A3 so both cells end up sharing a single cell style
# (mirrors the target workbook, which has exactly one extra cellXfs entry).
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"
$scratch.Font.Name = "Roboto"
$scratch.Font.Color = 3355443
$scratch.Copy()

$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("A2").Value = "00000001"
$ws.Range("A3").Value = "00000002"

$scratch.Clear()
$excel.CutCopyMode = $false

# New selection left by the author after the edit.
$ws.Range("B8").Select() | Out-Null

# Page setup now specifies paper size/orientation explicitly.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
